# Update cryptos list - price/volume refresh (GitHub Actions run, 2024-08-23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric (e.g. "576.18") while keeping
# the cell's underlying type as text, matching the source feed which stores
# every Price/Volume figure as a plain string (some values, like
# "60.720.70", use '.' as a thousands separator and are never ambiguous;
# others, like "576.18", would otherwise be auto-coerced to a number by
# Excel's input parser).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "60.720.70"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.641.66"
$ws.Range("E3").Value = "  +1.16%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "576.18"
$ws.Range("E5").Value = "  -0.55%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "143.86"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7 - USDC
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.34%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +0.61%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.05%  "

# Row 11 - Cardano
Set-TextValue $ws.Range("D11") "0.379"
$ws.Range("E11").Value = "  +2.08%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.85%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "3.110.47"
$ws.Range("E13").Value = "  +1.00%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "26.14"
$ws.Range("E14").Value = "  +11.60%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Range("D15") "60.701.81"
$ws.Range("E15").Value = "  -0.25%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.23%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.655.35"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +2.05%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "4.72"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "349.62"
$ws.Range("E20").Value = "  -0.42%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.85"
$ws.Range("E21").Value = "  -1.17%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.04%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +1.71%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "63.93"
$ws.Range("E24").Value = "  +1.04%  "

# Row 25 - Binance-PegBSC-USD
$ws.Range("E25").Value = "  -0.02%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +0.49%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +3.47%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  +9.28%  "

# Row 29 - PEPE
$ws.Range("E29").Value = "  +0.29%  "

# Row 30 - Aptos
$ws.Range("E30").Value = "  +6.60%  "

# Row 31 - USDe
$ws.Range("E31").Value = "  +0.09%  "

# Row 32 - Monero
Set-TextValue $ws.Range("D32") "163.36"
$ws.Range("E32").Value = "  +0.36%  "

# Row 33 - EthereumClassic
$ws.Range("E33").Value = "  +1.43%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  +7.14%  "

# Row 35 - Fetch.AI
$ws.Range("E35").Value = "  +3.44%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +6.84%  "

# Rows 37 & 38 swap ranking places: Bittensor overtakes Stacks
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D37") "339.26"
$ws.Range("E37").Value = "  +10.06%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D38") "1.66"
$ws.Range("E38").Value = "  +2.11%  "

# Row 39 - Filecoin
$ws.Range("E39").Value = "  +4.26%  "

# Row 40 - SuiNetwork
$ws.Range("E40").Value = "  +6.53%  "

# Row 42 - RenderToken
Set-TextValue $ws.Range("D42") "5.18"
$ws.Range("E42").Value = "  +2.43%  "

# Row 43 - Mantle
$ws.Range("E43").Value = "  +2.15%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "20.25"
$ws.Range("E44").Value = "  +1.53%  "

# Rows 45 & 46 swap ranking places: Hedera overtakes VeChain
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D45") "0.0562"
$ws.Range("E45").Value = "  +2.24%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0249"
$ws.Range("E46").Value = "  +2.41%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "132.85"
$ws.Range("E47").Value = "  -1.12%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +0.90%  "

# Row 49 - InjectiveProtocol
$ws.Range("E49").Value = "  +0.27%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  +0.38%  "

# Row 51 - Maker
Set-TextValue $ws.Range("D51") "2.086.79"
$ws.Range("E51").Value = "  +2.00%  "
